# update format_input and input_files table
#
# - Rewrites column A (original filename) on rows 2-6: "input_sbs/raw/..." -> "input/raw/..."
# - Rewrites the "snakemake filename" formula (col H) on rows 2-7 to branch on
#   tag (sbs vs phenotype) and include the well in the filename.
# - Adds a new row 7 for the phenotype ("c0-DAPI-GFP") channel.
# - Widens column H and splits out column C (cycle) to its own width.
# - Updates dimension / selection bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input_ascp")

# ---------------------------------------------------------------------------
# 1. Column A: fix the "original filename" source path prefix on rows 2-6
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "input/raw/10X_c4_B3_A594_Site-0.tif"
$ws.Range("A3").Value = "input/raw/10X_c4_B3_CY3_Site-0.tif"
$ws.Range("A4").Value = "input/raw/10X_c4_B3_CY5_Site-0.tif"
$ws.Range("A5").Value = "input/raw/10X_c4_B3_CY7_Site-0.tif"
$ws.Range("A6").Value = "input/raw/10X_c4_B3_DAPI_Site-0.tif"

# ---------------------------------------------------------------------------
# 2. Column H: new formula that branches sbs vs. other tags and inserts well
# ---------------------------------------------------------------------------
$formula2 = '=IF(F2="sbs","input/"&B2&"_c"&C2&"-SBS-"&C2&"/"&B2&"_c"&C2&"-SBS-"&C2&"_"&D2&"_Tile-"&E2&"."&F2&".tif","input/"&B2&"_"&C2&"/"&B2&"_"&C2&"_"&D2&"_Tile-"&E2&"."&F2&".tif")'
$ws.Range("H2").Formula = $formula2

$formula3 = '=IF(F3="sbs","input/"&B3&"_c"&C3&"-SBS-"&C3&"/"&B3&"_c"&C3&"-SBS-"&C3&"_"&D3&"_Tile-"&E3&"."&F3&".tif","input/"&B3&"_"&C3&"/"&B3&"_"&C3&"_"&D3&"_Tile-"&E3&"."&F3&".tif")'
$ws.Range("H3").Formula = $formula3

$formula4 = '=IF(F4="sbs","input/"&B4&"_c"&C4&"-SBS-"&C4&"/"&B4&"_c"&C4&"-SBS-"&C4&"_"&D4&"_Tile-"&E4&"."&F4&".tif","input/"&B4&"_"&C4&"/"&B4&"_"&C4&"_"&D4&"_Tile-"&E4&"."&F4&".tif")'
$ws.Range("H4").Formula = $formula4

$formula5 = '=IF(F5="sbs","input/"&B5&"_c"&C5&"-SBS-"&C5&"/"&B5&"_c"&C5&"-SBS-"&C5&"_"&D5&"_Tile-"&E5&"."&F5&".tif","input/"&B5&"_"&C5&"/"&B5&"_"&C5&"_"&D5&"_Tile-"&E5&"."&F5&".tif")'
$ws.Range("H5").Formula = $formula5

$formula6 = '=IF(F6="sbs","input/"&B6&"_c"&C6&"-SBS-"&C6&"/"&B6&"_c"&C6&"-SBS-"&C6&"_"&D6&"_Tile-"&E6&"."&F6&".tif","input/"&B6&"_"&C6&"/"&B6&"_"&C6&"_"&D6&"_Tile-"&E6&"."&F6&".tif")'
$ws.Range("H6").Formula = $formula6

# ---------------------------------------------------------------------------
# 3. Row 7: new phenotype ("c0-DAPI-GFP") entry
#    Seed formatting by duplicating row 6, then overwrite the values/formula.
# ---------------------------------------------------------------------------
$ws.Range("A6:H6").Copy($ws.Range("A7:H7"))

$ws.Range("C7").Value = "c0-DAPI-GFP"
$ws.Range("C7").ClearFormats()
$ws.Range("A7").Value = "input/raw/10X_c0-DAPI-GFP_B3_Site-0.tif"
$ws.Range("B7").Value = "10X"
$ws.Range("D7").Value = "B3"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "phenotype"
$ws.Range("G7").Value = "ALL"

$formula7 = '=IF(F7="sbs","input/"&B7&"_c"&C7&"-SBS-"&C7&"/"&B7&"_c"&C7&"-SBS-"&C7&"_"&D7&"_Tile-"&E7&"."&F7&".tif","input/"&B7&"_"&C7&"/"&B7&"_"&C7&"_"&D7&"_Tile-"&E7&"."&F7&".tif")'
$ws.Range("H7").Formula = $formula7

# ---------------------------------------------------------------------------
# 4. Column widths: split cycle (C) into its own width, widen snakemake
#    filename (H) to fit the longer formula result.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.5
$ws.Columns.Item(8).ColumnWidth = 43.333333333333336

# ---------------------------------------------------------------------------
# 5. Selection bookkeeping (cosmetic, matches author's saved cursor position)
# ---------------------------------------------------------------------------
$ws.Range("H12").Select()
